$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns L, M, N
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the look of the existing header row (bold, centered, bordered)
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# Row 2
$ws.Range("L2").Value = 91.70937389105484
$ws.Range("M2").Value = 221398
$ws.Range("N2").Value = 306.6454293628809

# Row 3
$ws.Range("L3").Value = 89.17110346557384
$ws.Range("M3").Value = 42155
$ws.Range("N3").Value = 390.3240740740741

# Row 4
$ws.Range("L4").Value = 89.84034764490647
$ws.Range("M4").Value = 173731
$ws.Range("N4").Value = 142.9884773662552

# Row 5
$ws.Range("L5").Value = 92.06346317438916
$ws.Range("M5").Value = 29915
$ws.Range("N5").Value = 178.0654761904762

# Row 6
$ws.Range("L6").Value = 19.49107025725424
$ws.Range("M6").Value = 2110
$ws.Range("N6").Value = 14.16107382550336

# Row 7
$ws.Range("L7").Value = 21.02856644875892
$ws.Range("M7").Value = 98
$ws.Range("N7").Value = 32.66666666666666
